$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $c = $ws.Range($CellRef)
    $c.NumberFormat = "@"
    $c.Value = $NewValue
    $c.Style = "Normal"
}

Set-TextValue 'D2' '68.069.84'
Set-TextValue 'E2' '  +0.34%  '
Set-TextValue 'D3' '3.247.17'
Set-TextValue 'E3' '  -0.02%  '
Set-TextValue 'E4' '  +0.01%  '
Set-TextValue 'D5' '582.05'
Set-TextValue 'D6' '184.71'
Set-TextValue 'E6' '  +0.97%  '
Set-TextValue 'E7' '  +0.00%  '
Set-TextValue 'E8' '  +0.93%  '
Set-TextValue 'E9' '  -3.20%  '
Set-TextValue 'E10' '  -0.89%  '
Set-TextValue 'E11' '  +0.33%  '
Set-TextValue 'D12' '3.808.87'
Set-TextValue 'E12' '  +0.00%  '
Set-TextValue 'D13' '0.137'
Set-TextValue 'E13' '  +0.17%  '
Set-TextValue 'D14' '27.90'
Set-TextValue 'E14' '  -2.84%  '
Set-TextValue 'D15' '68.073.12'
Set-TextValue 'E15' '  +0.34%  '
Set-TextValue 'E16' '  -0.80%  '
Set-TextValue 'D17' '3.238.63'
Set-TextValue 'E17' '  -0.44%  '
Set-TextValue 'E18' '  -0.30%  '
Set-TextValue 'D19' '13.49'
Set-TextValue 'E19' '  -0.49%  '
Set-TextValue 'D20' '396.51'
Set-TextValue 'E20' '  +4.50%  '
Set-TextValue 'E22' '  +0.25%  '
Set-TextValue 'D23' '71.35'
Set-TextValue 'E23' '  -0.02%  '
Set-TextValue 'E24' '  +0.74%  '
Set-TextValue 'E25' '  -0.49%  '
Set-TextValue 'D26' '0.187'
Set-TextValue 'E26' '  +2.59%  '
Set-TextValue 'D27' '9.66'
Set-TextValue 'E27' '  -2.41%  '
Set-TextValue 'E28' '  -0.06%  '
Set-TextValue 'E29' '  -0.34%  '
Set-TextValue 'E30' '  -1.05%  '
Set-TextValue 'D31' '22.83'
Set-TextValue 'E31' '  -0.03%  '
Set-TextValue 'D32' '7.03'
Set-TextValue 'E32' '  -0.36%  '
Set-TextValue 'D33' '1.27'
Set-TextValue 'E33' '  +0.48%  '
Set-TextValue 'D34' '0.999'
Set-TextValue 'E34' '  +0.06%  '
Set-TextValue 'D35' '161.90'
Set-TextValue 'E35' '  -0.18%  '
Set-TextValue 'E36' '  -5.06%  '
Set-TextValue 'D37' '1.91'
Set-TextValue 'E37' '  +3.10%  '
Set-TextValue 'B38' 'EnergySwap'
Set-TextValue 'C38' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D38' '26.74'
Set-TextValue 'E38' '  +1.12%  '
Set-TextValue 'B39' 'Mantle'
Set-TextValue 'C39' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D39' '0.814'
Set-TextValue 'E39' '  -2.97%  '
Set-TextValue 'D40' '4.59'
Set-TextValue 'E40' '  +0.39%  '
Set-TextValue 'E41' '  -2.76%  '
Set-TextValue 'E42' '  -3.87%  '
Set-TextValue 'D43' '41.23'
Set-TextValue 'E43' '  +0.18%  '
Set-TextValue 'D44' '0.0684'
Set-TextValue 'E44' '  -0.06%  '
Set-TextValue 'D45' '25.25'
Set-TextValue 'E45' '  -0.94%  '
Set-TextValue 'D46' '2.612.59'
Set-TextValue 'E46' '  -0.60%  '
Set-TextValue 'D47' '336.64'
Set-TextValue 'E47' '  -2.74%  '
Set-TextValue 'D48' '0.0280'
Set-TextValue 'E48' '  -1.28%  '
Set-TextValue 'E49' '  +2.37%  '
Set-TextValue 'E50' '  -0.98%  '
Set-TextValue 'B51' 'Arweave'
Set-TextValue 'C51' 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextValue 'D51' '31.14'
Set-TextValue 'E51' '  +2.78%  '
